$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.353.41"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.667.89"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.92%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5347"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2666"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06406"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.572"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "1.666.52"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").Value = "1.895.49"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5541"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "0.0₅8200"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "26.374.25"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.697"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.054"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1234"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.234"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.501"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05869"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.643"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.290"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.609"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9721"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.832"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.421"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5846"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01602"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.065.12"
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.845"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.011"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").Value = "1.806.83"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("E47").Value = "  -4.73%  "
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4387"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.017"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("E51").Value = "  +0.50%  "
